$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.699.71"
$ws.Range("E2").Value = "  -0.84%  "

$ws.Range("D3").Value = "1.901.05"
$ws.Range("E3").Value = "  -0.26%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9995"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.65%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "311.92"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.65%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9999"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.54%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4972"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.80%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3762"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.99%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07236"
$ws.Range("D9").Style = "Normal"

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.00"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.10%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.8908"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -4.42%  "

$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").Value = "1.940.96"
$ws.Range("E12").Value = "  +1.56%  "

$ws.Range("B13").Value = "TRON"
$ws.Range("C13").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07608"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.84%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.436"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.98%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "91.62"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.22%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.9995"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.73%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008682"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.02%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.9998"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.53%  "

$ws.Range("D19").Value = "27.740.89"
$ws.Range("E19").Value = "  -0.93%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.46"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.65%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.131"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.54%  "

$ws.Range("D22").Value = "2.125.56"
$ws.Range("E22").Value = "  -3.26%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.80"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.96%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.574"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.90%  "

$ws.Range("E25").Value = "  -1.78%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.840"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -4.02%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.190"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.48%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.25"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.15%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "114.59"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.38%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.842"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.35%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08902"
$ws.Range("D31").Style = "Normal"

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.185"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.36%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.770"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.40%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.225"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.81%  "

$ws.Range("E35").Value = "  +1.39%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.610"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.92%  "

$ws.Range("E37").Value = "  +0.67%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.052"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.63%  "

$ws.Range("E39").Value = "  -0.96%  "

$ws.Range("B40").Value = "Hedera"
$ws.Range("C40").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.05284"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.07%  "

$ws.Range("B41").Value = "TheSandbox"
$ws.Range("C41").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5490"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.03%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.744"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.15%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "113.66"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.54%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.437"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.42%  "

$ws.Range("E45").Value = "  -1.45%  "

$ws.Range("B46").Value = "Decentraland"
$ws.Range("C46").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4769"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.98%  "

$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.43"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.86%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.9993"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.61%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.614"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.45%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "66.76"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.73%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06007"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.24%  "
